$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.05970899999999999
$ws.Range("H2").Value = 0.179127
$ws.Range("I2").Value = 0.003688274646134975
$ws.Range("J2").Value = 0.003688274646134975
$ws.Range("M2").Value = 0.07271233333333334
$ws.Range("N2").Value = 0.218137
$ws.Range("O2").Value = 0.004171225362010892
$ws.Range("P2").Value = 0.004171225362010893
$ws.Range("Q2").Value = 0.004341580711
$ws.Range("R2").Value = 0.039074226399
$ws.Range("S2").Value = 0.00001538462474601996
$ws.Range("T2").Value = 0.00001538462474601996
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.05970899999999999
$ws.Range("H3").Value = 0.179127
$ws.Range("I3").Value = 0.003688274646134975
$ws.Range("J3").Value = 0.003688274646134975
$ws.Range("O3").Value = 0.5387060579248023
$ws.Range("P3").Value = 0.5387060579248023
$ws.Range("Q3").Value = 0.5607071368729999
$ws.Range("R3").Value = 5.046364231856999
$ws.Range("S3").Value = 0.001986895895163367
$ws.Range("T3").Value = 0.001986895895163368
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.05970899999999999
$ws.Range("H4").Value = 0.179127
$ws.Range("I4").Value = 0.003688274646134975
$ws.Range("J4").Value = 0.003688274646134975
$ws.Range("O4").Value = 0.4571227167131868
$ws.Range("P4").Value = 0.4571227167131868
$ws.Range("Q4").Value = 0.475791883008
$ws.Range("R4").Value = 4.282126947071999
$ws.Range("S4").Value = 0.001685994126225587
$ws.Range("T4").Value = 0.001685994126225588
$ws.Range("I5").Value = 0.8850509663933519
$ws.Range("J5").Value = 0.885050966393352
$ws.Range("M5").Value = 0.07271233333333334
$ws.Range("N5").Value = 0.218137
$ws.Range("O5").Value = 0.004171225362010892
$ws.Range("P5").Value = 0.004171225362010893
$ws.Range("Q5").Value = 1.041820518429111
$ws.Range("R5").Value = 9.376384665862
$ws.Range("S5").Value = 0.0036917470376922
$ws.Range("T5").Value = 0.0036917470376922
$ws.Range("I6").Value = 0.8850509663933519
$ws.Range("J6").Value = 0.885050966393352
$ws.Range("O6").Value = 0.5387060579248023
$ws.Range("P6").Value = 0.5387060579248023
$ws.Range("S6").Value = 0.4767823171682993
$ws.Range("T6").Value = 0.4767823171682993
$ws.Range("I7").Value = 0.8850509663933519
$ws.Range("J7").Value = 0.885050966393352
$ws.Range("O7").Value = 0.4571227167131868
$ws.Range("P7").Value = 0.4571227167131868
$ws.Range("S7").Value = 0.4045769021873604
$ws.Range("T7").Value = 0.4045769021873604
$ws.Range("I8").Value = 0.111260758960513
$ws.Range("J8").Value = 0.111260758960513
$ws.Range("M8").Value = 0.07271233333333334
$ws.Range("N8").Value = 0.218137
$ws.Range("O8").Value = 0.004171225362010892
$ws.Range("P8").Value = 0.004171225362010893
$ws.Range("Q8").Value = 0.1309684368273334
$ws.Range("R8").Value = 1.178715931446
$ws.Range("S8").Value = 0.0004640936995726725
$ws.Range("T8").Value = 0.0004640936995726727
$ws.Range("I9").Value = 0.111260758960513
$ws.Range("J9").Value = 0.111260758960513
$ws.Range("O9").Value = 0.5387060579248023
$ws.Range("P9").Value = 0.5387060579248023
$ws.Range("S9").Value = 0.05993684486133959
$ws.Range("T9").Value = 0.0599368448613396
$ws.Range("I10").Value = 0.111260758960513
$ws.Range("J10").Value = 0.111260758960513
$ws.Range("O10").Value = 0.4571227167131868
$ws.Range("P10").Value = 0.4571227167131868
$ws.Range("S10").Value = 0.05085982039960075
$ws.Range("T10").Value = 0.05085982039960076
